$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow writes, then restore protection at the end.
$ws.Unprotect("")

# --- Text corrections in shared strings ---
# Company name correction: "D.R. Horton Inc" -> "D R Horton Inc" (row 55, DHI)
$ws.Range("B55").Value = "D R Horton Inc"

# Update the "as of" date in the disclaimer footer (row 80)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# --- Updated Weight (D) and Percent Change (E) values for rows 2-77 ---
$ws.Range("D2").Value = 0.06303203901975518
$ws.Range("E2").Value = -0.009258532757944415
$ws.Range("D3").Value = 0.03795063583417754
$ws.Range("E3").Value = 0.01473517639393096
$ws.Range("D4").Value = 0.03214246582418404
$ws.Range("E4").Value = -0.01196856739875063
$ws.Range("D5").Value = 0.02939516902736379
$ws.Range("E5").Value = -0.002203225522164454
$ws.Range("D6").Value = 0.02682862318777295
$ws.Range("E6").Value = 0.004626093979055224
$ws.Range("D7").Value = 0.02575024060385971
$ws.Range("E7").Value = 0.004024144869215318
$ws.Range("D8").Value = 0.1904852107349112
$ws.Range("E8").Value = 0.002369668246445578
$ws.Range("D9").Value = 0.02472084329565163
$ws.Range("E9").Value = 0.0009987075549287017
$ws.Range("D10").Value = 0.02275154100985069
$ws.Range("E10").Value = -0.002028838489964491
$ws.Range("D11").Value = 0.02205985949664843
$ws.Range("E11").Value = -0.002903811252268684
$ws.Range("D12").Value = 0.01909012715702275
$ws.Range("E12").Value = -0.02084052964881955
$ws.Range("D13").Value = 0.02028461426460256
$ws.Range("E13").Value = 0.008970727101038856
$ws.Range("D14").Value = 0.01723394631720121
$ws.Range("E14").Value = 0.0007561436672967048
$ws.Range("D15").Value = 0.01613640923445454
$ws.Range("E15").Value = -0.006947079599521344
$ws.Range("D16").Value = 0.01465413363210688
$ws.Range("E16").Value = -0.01114253878460614
$ws.Range("D17").Value = 0.01423838680496817
$ws.Range("E17").Value = -0.004515481651376385
$ws.Range("D18").Value = 0.01438451836061747
$ws.Range("E18").Value = -0.00286513404734301
$ws.Range("D19").Value = 0.01364106474913386
$ws.Range("E19").Value = -0.001519275811863063
$ws.Range("D20").Value = 0.01335759386680798
$ws.Range("E20").Value = 0.0233667928254071
$ws.Range("D21").Value = 0.01252798426159613
$ws.Range("E21").Value = -0.02698511166253104
$ws.Range("D22").Value = 0.0132529898926459
$ws.Range("E22").Value = 0.003950871768444442
$ws.Range("D23").Value = 0.01142707159094029
$ws.Range("E23").Value = -0.001676238630430649
$ws.Range("D24").Value = 0.01286805511793392
$ws.Range("E24").Value = -0.002855051244509577
$ws.Range("D25").Value = 0.01143252748302604
$ws.Range("E25").Value = -0.01118561320025957
$ws.Range("D26").Value = 0.008759258114074343
$ws.Range("E26").Value = 0.008657465495608641
$ws.Range("D27").Value = 0.009467189550467266
$ws.Range("E27").Value = -0.01810975306389817
$ws.Range("D28").Value = 0.01008060452540343
$ws.Range("E28").Value = -0.002032520325203402
$ws.Range("D29").Value = 0.009987658104834653
$ws.Range("E29").Value = 0.002137894174238486
$ws.Range("D30").Value = 0.009778214650377308
$ws.Range("E30").Value = -0.002729608220937596
$ws.Range("D31").Value = 0.008472921907127883
$ws.Range("E31").Value = -0.0003613369467029282
$ws.Range("D32").Value = 0.01059494992030767
$ws.Range("E32").Value = -0.05495910020449912
$ws.Range("D33").Value = 0.009444934220880065
$ws.Range("E33").Value = -0.01277900834895207
$ws.Range("D34").Value = 0.009022475468945209
$ws.Range("E34").Value = -0.001644436323771092
$ws.Range("D35").Value = 0.009361368794617004
$ws.Range("E35").Value = -0.003983228511530434
$ws.Range("D36").Value = 0.008358819185593494
$ws.Range("E36").Value = -0.01327022229735453
$ws.Range("D37").Value = 0.008723147173650664
$ws.Range("E37").Value = -0.00492710583153344
$ws.Range("D38").Value = 0.006944369349606834
$ws.Range("E38").Value = -0.02189100281479961
$ws.Range("D39").Value = 0.008796625087208413
$ws.Range("E39").Value = 0.002248875562218755
$ws.Range("D40").Value = 0.008163937860372221
$ws.Range("E40").Value = 0.01224079656526911
$ws.Range("D41").Value = 0.006933614569524129
$ws.Range("E41").Value = -0.003804175535528254
$ws.Range("D42").Value = 0.007151928754998594
$ws.Range("E42").Value = -0.002689204763734132
$ws.Range("D43").Value = 0.008050934167459422
$ws.Range("E43").Value = 0.002408416782862188
$ws.Range("D44").Value = 0.007469548026638177
$ws.Range("E44").Value = 0.003531229309203354
$ws.Range("D45").Value = 0.007207586704477699
$ws.Range("E45").Value = -0.01339664974840438
$ws.Range("D46").Value = 0.008033585215647177
$ws.Range("E46").Value = 0.006097560975609762
$ws.Range("D47").Value = 0.007375110067225815
$ws.Range("E47").Value = 0.0201813769319199
$ws.Range("D48").Value = 0.007189688238354659
$ws.Range("E48").Value = -0.001255650426921107
$ws.Range("D49").Value = 0.006547070502901956
$ws.Range("E49").Value = -0.001978417266186971
$ws.Range("D50").Value = 0.007282085144612639
$ws.Range("E50").Value = 0.002727380528874779
$ws.Range("D51").Value = 0.006655442575195195
$ws.Range("E51").Value = 0.01056846798497291
$ws.Range("D52").Value = 0.006726604678443164
$ws.Range("E52").Value = -0.008647752867996239
$ws.Range("D53").Value = 0.005453144514197655
$ws.Range("E53").Value = 0.009501187648455867
$ws.Range("D54").Value = 0.006161154452634978
$ws.Range("E54").Value = -0.003465674532388774
$ws.Range("D55").Value = 0.005669221391406739
$ws.Range("E55").Value = -0.006646588430782141
$ws.Range("D56").Value = 0.005701073595921756
$ws.Range("E56").Value = -0.001869236091747606
$ws.Range("D57").Value = 0.006794508946848562
$ws.Range("E57").Value = -0.0004159349292910397
$ws.Range("D58").Value = 0.005512727565896726
$ws.Range("E58").Value = -0.003417634996582319
$ws.Range("D59").Value = 0.005409104867289644
$ws.Range("E59").Value = 0.01123301985370961
$ws.Range("D60").Value = 0.004959445156970431
$ws.Range("E60").Value = -0.003039128783082101
$ws.Range("D61").Value = 0.004925100512545736
$ws.Range("E61").Value = 0.01236083106864205
$ws.Range("D62").Value = 0.005062989353533114
$ws.Range("E62").Value = 0.006977285060857152
$ws.Range("D63").Value = 0.004213361726999446
$ws.Range("E63").Value = -0.004098971530779649
$ws.Range("D64").Value = 0.004074570112501476
$ws.Range("E64").Value = 0.008939580764488309
$ws.Range("D65").Value = 0.003805465104300662
$ws.Range("E65").Value = -0.005322221305388175
$ws.Range("D66").Value = 0.003784112548224051
$ws.Range("E66").Value = 0.003360716952949927
$ws.Range("D67").Value = 0.003869365768441696
$ws.Range("E67").Value = -0.004990870359099397
$ws.Range("D68").Value = 0.003643633139771497
$ws.Range("E68").Value = -0.0005493972788676027
$ws.Range("D69").Value = 0.00355418006017849
$ws.Range("E69").Value = -0.005057979017117686
$ws.Range("D70").Value = 0.002992968944767865
$ws.Range("E70").Value = 0.002885170225043199
$ws.Range("D71").Value = 0.002900140277265688
$ws.Range("E71").Value = -0.007836290551788472
$ws.Range("D72").Value = 0.002236209236758574
$ws.Range("E72").Value = -0.005441269395492609
$ws.Range("D73").Value = 0.001936527682264231
$ws.Range("E73").Value = -0.008979062366986201
$ws.Range("D74").Value = 0.001909797736146268
$ws.Range("E74").Value = -0.007851035843472576
$ws.Range("D75").Value = 0.001494678925362749
$ws.Range("E75").Value = 0.0117647058823529
$ws.Range("D76").Value = 0.001710795053594033
$ws.Range("E76").Value = 0.03280870004129777
$ws.Range("E77").Value = -0.001506670112644715

# Restore sheet protection (best effort match of original settings)
$ws.Protect()
